$wb = $excel.ActiveWorkbook

# --- ATEO4 ---
$ws = $wb.Worksheets.Item("ATEO4")
$ws.Rows.Item(40).Delete()
$ws.Range("A1").Value = 'Тип АС'
$ws.Range("C1").Value = 'Широкополосная'
$ws.Range("A2").Value = 'Пиковая мощность'
$ws.Range("C2").Value = '140 W'
$ws.Range("A3").Value = 'Program power handling'
$ws.Range("C3").Value = '70 W'
$ws.Range("A4").Value = 'Продолжительная мощность'
$ws.Range("C4").Value = '35 W'
$ws.Range("A5").Value = 'Impedance'
$ws.Range("C5").Value = '8 Ω (ATEO4)'
$ws.Range("A6").Value = '-'
$ws.Range("C6").Value = '16 Ω (ATEO4D)'
$ws.Range("A7").Value = 'Incline angle'
$ws.Range("C7").Value = '-'
$ws.Range("A8").Value = '          Left'
$ws.Range("C8").Value = '30°'
$ws.Range("A9").Value = '          Right'
$ws.Range("C9").Value = '30°'
$ws.Range("A10").Value = '          Bottom'
$ws.Range("C10").Value = '30°'
$ws.Range("A11").Value = '          Top'
$ws.Range("C11").Value = '5°'
$ws.Range("A12").Value = 'Sensitivity (1W/1m)'
$ws.Range("C12").Value = '86 dB'
$ws.Range("A13").Value = 'Sound Pressure (Max. W/1m)'
$ws.Range("C13").Value = '-'
$ws.Range("A14").Value = '          @ 8 Ω'
$ws.Range("C14").Value = '101 dB'
$ws.Range("A15").Value = '          @ 100 V'
$ws.Range("C15").Value = '100 dB'
$ws.Range("A16").Value = 'Frequency'
$ws.Range("C16").Value = '-'
$ws.Range("A17").Value = '          Отклик (± 3 дБ)'
$ws.Range("C17").Value = '100 Hz - 20 kHz'
$ws.Range("A18").Value = '          Range (-10 dB)'
$ws.Range("C18").Value = '65 Hz - 20 kHz'
$ws.Range("A19").Value = 'Dispersion'
$ws.Range("C19").Value = '-'
$ws.Range("A20").Value = '          Horizontal'
$ws.Range("C20").Value = '130°'
$ws.Range("A21").Value = '          Vertical'
$ws.Range("C21").Value = '130°'
$ws.Range("A22").Value = 'Connectors'
$ws.Range("C22").Value = '6-pin Custom Terminal Block'
$ws.Range("A23").Value = 'Drivers'
$ws.Range("C23").Value = '1” Dome tweeter'
$ws.Range("A24").Value = '-'
$ws.Range("C24").Value = '4” MF / LF Woofer'
$ws.Range("A25").Value = 'Line Transformer Tappings'
$ws.Range("C25").Value = '-'
$cell = $ws.Range("A26")
$cell.Formula = '="          1"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("C26").Value = '208 Ω - 100 В / Нет - 70 В / 24 Вт'
$cell = $ws.Range("A27")
$cell.Formula = '="          2"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("C27").Value = '417 Ω - 100 V / 24 W - 70 V / 12 W'
$cell = $ws.Range("A28")
$cell.Formula = '="          3"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("C28").Value = '833 Ω - 100 V / 12 W - 70 V / 6 W'
$cell = $ws.Range("A29")
$cell.Formula = '="          4"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("C29").Value = '1667 Ω - 100 V / 6 W - 70 V / 3 W'
$ws.Range("A30").Value = 'Dimensions'
$ws.Range("C30").Value = '136 x 244 x 153 mm (W x H x D)'
$ws.Range("A31").Value = 'Weight'
$ws.Range("C31").Value = '2.050 kg'
$ws.Range("A32").Value = 'Operating temperature'
$ws.Range("C32").Value = '-20 °C ~ 60 °C'
$ws.Range("A33").Value = 'Mounting'
$ws.Range("C33").Value = 'Clevermount™'
$ws.Range("A34").Value = 'Accessories'
$ws.Range("C34").Value = '-'
$ws.Range("A35").Value = '          Included'
$ws.Range("C35").Value = 'Allen-key wrench'
$ws.Range("A36").Value = 'Construction'
$ws.Range("C36").Value = 'ABS'
$ws.Range("A37").Value = 'Front finish'
$ws.Range("C37").Value = 'Steel grill'
$ws.Range("A38").Value = 'Colours'
$ws.Range("C38").Value = 'Black (RAL9005) (ATEO4/B, ATEO4D/B)'
$ws.Range("A39").Value = '-'
$ws.Range("C39").Value = 'White (RAL9003) (ATEO4/W, ATEO4D/W)'
$excel.CutCopyMode = $false

# --- ALTI6 ---
$ws = $wb.Worksheets.Item("ALTI6")
$ws.Rows.Item(29).Delete()
$ws.Range("A1").Value = 'Тип АС'
$ws.Range("C1").Value = '2-way coaxial'
$ws.Range("A2").Value = 'Пиковая мощность'
$ws.Range("C2").Value = '240 W'
$ws.Range("A3").Value = 'Program power handling'
$ws.Range("C3").Value = '120 W'
$ws.Range("A4").Value = 'Продолжительная мощность'
$ws.Range("C4").Value = '60 W'
$ws.Range("A5").Value = 'Impedance'
$ws.Range("C5").Value = '16 Ω'
$ws.Range("A6").Value = 'Line Transformer Tappings'
$ws.Range("C6").Value = '-'
$cell = $ws.Range("A7")
$cell.Formula = '="          1"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("C7").Value = '83 Ω - 100 V / Not used - 70 V / 60 W'
$cell = $ws.Range("A8")
$cell.Formula = '="          2"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("C8").Value = '167 Ω - 100 V / 60 W - 70 V / 30 W'
$cell = $ws.Range("A9")
$cell.Formula = '="          3"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("C9").Value = '333 Ω - 100 V / 30 W - 70 V / 15 W'
$cell = $ws.Range("A10")
$cell.Formula = '="          4"'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("C10").Value = '667 Ω - 100 V / 15 W - 70 V / 7.5 W'
$ws.Range("A11").Value = 'Dispersion'
$ws.Range("C11").Value = '-'
$ws.Range("A12").Value = '          Conical'
$ws.Range("C12").Value = '115° (average 500 Hz to 5 kHz @ -6 dB)'
$ws.Range("A13").Value = 'Connectors'
$ws.Range("C13").Value = '4-pin Euro Terminal Block'
$ws.Range("A14").Value = 'Sensitivity (1W/1m)'
$ws.Range("C14").Value = '83 dB'
$ws.Range("A15").Value = 'Sound Pressure (Max. W/1m)'
$ws.Range("C15").Value = '-'
$ws.Range("A16").Value = '          @ 16 Ω'
$ws.Range("C16").Value = '101 dB'
$ws.Range("A17").Value = 'Frequency'
$ws.Range("C17").Value = '-'
$ws.Range("A18").Value = '          Range (-10 dB)'
$ws.Range("C18").Value = '61.5 Hz - 20 kHz'
$ws.Range("A19").Value = '          Отклик (± 3 дБ)'
$ws.Range("C19").Value = '75 Hz - 17 kHz'
$ws.Range("A20").Value = 'Dimensions'
$ws.Range("C20").Value = '230 x 302 mm (Ø x H)'
$ws.Range("A21").Value = 'Connection cable length'
$ws.Range("C21").Value = '3.5 m'
$ws.Range("A22").Value = 'Construction'
$ws.Range("C22").Value = 'Polypropylene'
$ws.Range("A23").Value = 'Front finish'
$ws.Range("C23").Value = 'Aluminium grill'
$ws.Range("A24").Value = 'Mounting & handling'
$ws.Range("C24").Value = 'Dual snap hook and dual Gripple™ hanger'
$ws.Range("A25").Value = 'Accessories'
$ws.Range("C25").Value = '-'
$ws.Range("A26").Value = '          Included'
$ws.Range("C26").Value = 'Connection cable with dual steel core, snap hook & open ends – 3.5 meter'
$ws.Range("A27").Value = '          -'
$ws.Range("C27").Value = '2 Gripple™ hangers'
$ws.Range("A28").Value = 'Weight'
$ws.Range("C28").Value = '3.39 kg'
$excel.CutCopyMode = $false

# --- XENO6 ---
$ws = $wb.Worksheets.Item("XENO6")
$ws.Rows.Item(32).Delete()
$ws.Range("A1").Value = 'Тип АС'
$ws.Range("C1").Value = '2-way'
$ws.Range("A2").Value = 'Пиковая мощность'
$ws.Range("C2").Value = '320 W'
$ws.Range("A3").Value = 'Program power handling'
$ws.Range("C3").Value = '160 W'
$ws.Range("A4").Value = 'Продолжительная мощность'
$ws.Range("C4").Value = '80 W'
$ws.Range("A5").Value = 'Impedance'
$ws.Range("C5").Value = '8 Ω'
$ws.Range("A6").Value = 'Sensitivity (1W/1m)'
$ws.Range("C6").Value = '88 dB'
$ws.Range("A7").Value = 'Sound Pressure (Max. W/1m)'
$ws.Range("C7").Value = '110 dB'
$ws.Range("A8").Value = 'Frequency'
$ws.Range("C8").Value = '-'
$ws.Range("A9").Value = '          Отклик (± 3 дБ)'
$ws.Range("C9").Value = '65 Hz - 18 kHz'
$ws.Range("A10").Value = '          Range (-10 dB)'
$ws.Range("C10").Value = '55 Hz - 20 kHz'
$ws.Range("A11").Value = 'Crossover'
$ws.Range("C11").Value = '-'
$ws.Range("A12").Value = '          Frequency'
$ws.Range("C12").Value = '2.5 kHz'
$ws.Range("A13").Value = '          Type'
$ws.Range("C13").Value = 'Passive built-in'
$ws.Range("A14").Value = 'Dispersion'
$ws.Range("C14").Value = '-'
$ws.Range("A15").Value = '          Horizontal'
$ws.Range("C15").Value = '120°'
$ws.Range("A16").Value = '          Vertical'
$ws.Range("C16").Value = '120°'
$ws.Range("A17").Value = 'Connectors'
$ws.Range("C17").Value = '4-pin Euro Terminal Block (Pitch - 5.08 mm)'
$ws.Range("A18").Value = 'Drivers'
$ws.Range("C18").Value = '1” Dome tweeter'
$ws.Range("A19").Value = '-'
$ws.Range("C19").Value = '6” MF / LF Woofer'
$ws.Range("A20").Value = 'Dimensions'
$ws.Range("C20").Value = '210 x 345 x 220 mm (W x H x D)'
$ws.Range("A21").Value = 'Weight'
$ws.Range("C21").Value = '5.600 kg'
$ws.Range("A22").Value = 'Operating temperature'
$ws.Range("C22").Value = '-20 °C ~ 60 °C'
$ws.Range("A23").Value = 'Construction'
$ws.Range("C23").Value = 'Medium Density Fibreboard with structured coating'
$ws.Range("A24").Value = 'Front finish'
$ws.Range("C24").Value = 'Steel grill'
$ws.Range("A25").Value = 'Mounting & handling'
$ws.Range("C25").Value = '2-way Revolving mounting bracket'
$ws.Range("A26").Value = 'Colours'
$ws.Range("C26").Value = 'Black (RAL9004) (XENO6/B)'
$ws.Range("A27").Value = '-'
$ws.Range("C27").Value = 'White (RAL9003) (XENO6/W)'
$ws.Range("A28").Value = 'Accessories'
$ws.Range("C28").Value = '-'
$ws.Range("A29").Value = '          Included'
$ws.Range("C29").Value = '4-pin Euro Terminal Block (Pitch - 5.08 mm)'
$ws.Range("A30").Value = '          -'
$ws.Range("C30").Value = '2-way Revolving mounting bracket'
$ws.Range("A31").Value = '          Optional'
$ws.Range("C31").Value = 'WBP100 Wall bracket mounting plate'
$excel.CutCopyMode = $false

